$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "editUserExp1"
$ws.Range("B16").Value = "//tr[td[@class='sorting_1' ]/label[contains(text(),'"

$ws.Range("A17").Value = "editUserExp2"
$ws.Range("B17").Formula = "'') ]]/td/button[@id='edit']"

$ws.Range("B18").Value = "user_save_btn"
$ws.Range("A18").Value = "editUser"

$ws.Range("A18").Select() | Out-Null
